$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update Version, Status, Date, Contact ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "0.4.0-snapshot-1"
$meta.Range("B6").Value = "draft"
$meta.Range("B8").Value = "2024-05-23T12:16:26+00:00"
$meta.Range("B10").Value = "ANS (https://esante.gouv.fr)"

# --- Sheet "Elements": swap the two Mapping columns (AK <-> AL) ---
$els = $wb.Worksheets.Item("Elements")

# swap header text (row 1)
$akHeader = $els.Range("AK1").Value()
$alHeader = $els.Range("AL1").Value()
$els.Range("AK1").Value = $alHeader
$els.Range("AL1").Value = $akHeader

# swap data rows 2..6 (only rows whose AK/AL actually differ need touching;
# rows where both are already blank would otherwise be turned from an
# explicit empty-string cell into a fully-cleared cell)
for ($r = 2; $r -le 6; $r++) {
    $akCell = $els.Cells.Item($r, 37)
    $alCell = $els.Cells.Item($r, 38)
    $akVal = $akCell.Value()
    $alVal = $alCell.Value()
    if ($akVal -ne $alVal) {
        $akCell.Value = $alVal
        $alCell.Value = $akVal
    }
}

# swap column widths (AK was 24.98046875, AL was 77.30859375 - stored OOXML width)
$els.Columns.Item(37).ColumnWidth = 76.47526041666667
$els.Columns.Item(38).ColumnWidth = 24.147135416666668
